$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$refStyle = $ws.Range("B2")

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $refStyle.Style
}

Set-TextValue "D2" '43.786.91'
Set-TextValue "E2" '  -0.58%  '
Set-TextValue "D3" '2.343.61'
Set-TextValue "E3" '  -0.31%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '239.08'
Set-TextValue "E5" '  +0.08%  '
Set-TextValue "D6" '0.664'
Set-TextValue "E6" '  -1.91%  '
Set-TextValue "D7" '72.75'
Set-TextValue "E7" '  -1.62%  '
Set-TextValue "E8" '  -0.07%  '
Set-TextValue "D9" '0.595'
Set-TextValue "E9" '  +0.95%  '
Set-TextValue "E10" '  +0.45%  '
Set-TextValue "D11" '60.77'
Set-TextValue "E11" '  +6.31%  '
Set-TextValue "D12" '33.31'
Set-TextValue "E12" '  +3.40%  '
Set-TextValue "E13" '  +0.37%  '
Set-TextValue "D14" '7.18'
Set-TextValue "E14" '  -0.02%  '
Set-TextValue "D15" '16.07'
Set-TextValue "E15" '  -2.83%  '
Set-TextValue "D16" '0.897'
Set-TextValue "E16" '  -0.10%  '
Set-TextValue "D17" '2.336.57'
Set-TextValue "E17" '  -0.86%  '
Set-TextValue "D18" '43.738.76'
Set-TextValue "E18" '  -0.51%  '
Set-TextValue "E19" '  +0.00%  '
Set-TextValue "D20" '77.65'
Set-TextValue "E20" '  +1.27%  '
Set-TextValue "D21" '6.48'
Set-TextValue "E21" '  -3.04%  '
Set-TextValue "D22" '251.74'
Set-TextValue "E22" '  -1.65%  '
Set-TextValue "D23" '3.79'
Set-TextValue "E23" '  +2.47%  '
Set-TextValue "D24" '0.999'
Set-TextValue "E24" '  -0.07%  '
Set-TextValue "E25" '  -4.81%  '
Set-TextValue "D26" '2.48'
Set-TextValue "E26" '  -0.41%  '
Set-TextValue "D27" '10.37'
Set-TextValue "E27" '  -2.79%  '
Set-TextValue "D28" '2.26'
Set-TextValue "E28" '  +0.95%  '
Set-TextValue "D29" '175.51'
Set-TextValue "E29" '  +0.18%  '
Set-TextValue "D30" '22.14'
Set-TextValue "E30" '  -2.45%  '
Set-TextValue "E31" '  +0.21%  '
Set-TextValue "E32" '  -2.40%  '
Set-TextValue "D33" '0.0737'
Set-TextValue "E33" '  -2.66%  '
Set-TextValue "D34" '5.02'
Set-TextValue "E34" '  -4.71%  '
Set-TextValue "D35" '5.31'
Set-TextValue "E35" '  -0.38%  '
Set-TextValue "D36" '3.74'
Set-TextValue "E36" '  +0.52%  '
Set-TextValue "B37" 'THORChain'
Set-TextValue "C37" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D37" '6.40'
Set-TextValue "E37" '  +1.15%  '
Set-TextValue "B38" 'LidoDAOToken'
Set-TextValue "C38" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D38" '2.38'
Set-TextValue "E38" '  +1.22%  '
Set-TextValue "D39" '0.0270'
Set-TextValue "E39" '  -3.65%  '
Set-TextValue "E40" '  +15.77%  '
Set-TextValue "D41" '65.55'
Set-TextValue "E41" '  +14.35%  '
Set-TextValue "D42" '19.67'
Set-TextValue "E42" '  +2.97%  '
Set-TextValue "D43" '9.06'
Set-TextValue "E43" '  +0.04%  '
Set-TextValue "E44" '  -3.37%  '
Set-TextValue "E45" '  -2.61%  '
Set-TextValue "E46" '  -0.13%  '
Set-TextValue "E47" '  -2.09%  '
Set-TextValue "D48" '2.41'
Set-TextValue "E48" '  -2.49%  '
Set-TextValue "D49" '1.14'
Set-TextValue "E49" '  -2.12%  '
Set-TextValue "D50" '97.25'
Set-TextValue "E50" '  -2.68%  '
Set-TextValue "E51" '  +1.88%  '
